# Updating the Staging testdata
# The ExpectedFilenames column (G) is expanded so that each Study_Type
# (Clinical, Economic, Quality of Life, Real-world Evidence) now expects
# three report files - a "Standard" excel report, a plain excel report and
# a word report - instead of just one filename each. The "ExcelReport-"
# filenames also drop the space around the first hyphen after
# NewImportLogic_1. Column C / F keep the same values, they are simply
# re-pointed because some now-unused shared strings were removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (Study_Types) values are unchanged, just re-asserting them ---
$ws.Range("C2").Value  = "Clinical"
$ws.Range("C3").Value  = "Economic"
$ws.Range("C4").Value  = "Quality of life"
$ws.Range("C5").Value  = "Real-world Evidence"
$ws.Range("C7").Value  = "Clinical"
$ws.Range("C8").Value  = "Economic"
$ws.Range("C9").Value  = "Quality of life"
$ws.Range("C10").Value = "Real-world Evidence"
$ws.Range("C12").Value = "Clinical"

# --- Column F (Invalid_Files) values are unchanged, just re-asserting them ---
$ws.Range("F1").Value  = "Invalid_Files"
$ws.Range("F2").Value  = "\Testdata\Templates\ManageQAData\1stUpload\InvalidDocumentFormat_1.docx"
$ws.Range("F3").Value  = "\Testdata\Templates\ManageQAData\1stUpload\InvalidDocumentFormat_2.html"
$ws.Range("F4").Value  = "\Testdata\Templates\ManageQAData\1stUpload\InvalidDocumentFormat_3.txt"
$ws.Range("F5").Value  = "\Testdata\Templates\ManageQAData\1stUpload\InvalidDocumentFormat_4.PNG"
$ws.Range("F7").Value  = "\Testdata\Templates\ManageQAData\1stUpload\InvalidDocumentFormat_1.docx"
$ws.Range("F8").Value  = "\Testdata\Templates\ManageQAData\1stUpload\InvalidDocumentFormat_2.html"
$ws.Range("F9").Value  = "\Testdata\Templates\ManageQAData\1stUpload\InvalidDocumentFormat_3.txt"
$ws.Range("F10").Value = "\Testdata\Templates\ManageQAData\1stUpload\InvalidDocumentFormat_4.PNG"

# --- Column G (ExpectedFilenames) now lists Standard/Excel/Word report names
#     for every study type, one after another, continuing straight down the
#     column (rows 2-13) ---
$ws.Range("G1").Value  = "ExpectedFilenames"

$ws.Range("G2").Value  = "StandardExcelReport-NewImportLogic_1 - Test_Automation_1-Clinical-2023_"
$ws.Range("G3").Value  = "ExcelReport-NewImportLogic_1-Test_Automation_1-Clinical-"
$ws.Range("G4").Value  = "WordReport-NewImportLogic_1 - Test_Automation_1-Clinical-"

$ws.Range("G5").Value  = "StandardExcelReport-NewImportLogic_1 - Test_Automation_1-Economic-2023_"
$ws.Range("G6").Value  = "ExcelReport-NewImportLogic_1-Test_Automation_1-Economic-"
$ws.Range("G7").Value  = "WordReport-NewImportLogic_1 - Test_Automation_1-Economic-"

$ws.Range("G8").Value  = "StandardExcelReport-NewImportLogic_1 - Test_Automation_1-Quality of Life-2023_"
$ws.Range("G9").Value  = "ExcelReport-NewImportLogic_1-Test_Automation_1-Quality of Life-"
$ws.Range("G10").Value = "WordReport-NewImportLogic_1 - Test_Automation_1-Quality of Life-"

$ws.Range("G11").Value = "StandardExcelReport-NewImportLogic_1 - Test_Automation_1-Real-world Evidence-2023_"
$ws.Range("G12").Value = "ExcelReport-NewImportLogic_1-Test_Automation_1-Real-world Evidence-"
$ws.Range("G13").Value = "WordReport-NewImportLogic_1 - Test_Automation_1-Real-world Evidence-"

# New column G needs a best-fit custom width, like the other text columns.
$ws.Columns.Item(7).ColumnWidth = 75.6

# Move / collapse the active selection onto the new last filled cell.
$ws.Range("G13").Select()
